$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The account-statement data table (rows 16-32) is being refreshed with
# a new set of workers/periods ("parte 1 de nuevos estado de cuenta").
# The new table only has 14 data rows instead of 17, so first remove 3
# rows from inside the table; everything below (the blank spacer rows
# and the signature footer) will naturally shift up to follow.
# ---------------------------------------------------------------------
$ws.Range("A29:A31").EntireRow.Delete() | Out-Null

# New data for rows 16-29 : TipoDoc | NumDoc | Nombre | Periodo | ValorMora | SalarioBasico
$data = @(
    ,@("CC", "73588295",   "JESUS MARIA BOSSIO HERRERA",       "1612", 15600, 900000)
    ,@("CC", "20035381",   "EISTON RAFAEL CABARCAS BELTRAN",   "1612", 15600, 900000)
    ,@("CC", "1047464278", "ANYELO SIMANCAS BARRIOS",          "1712", 2395,  898174)
    ,@("CC", "20144861",   "ARTURO RAFAEL MARIO PORTO",        "1806", 10660, 1230000)
    ,@("CC", "53105458",   "VIVIANA ROJAS ECHEVERRY",          "1809", 20000, 6000000)
    ,@("CC", "1075240833", "ADRIANA DEL PILAR REYES HUEPENDO", "1809", 6000,  1800000)
    ,@("CC", "20187337",   "YAIR ENRIQUEZ CASTILLO AMAYA",     "1902", 9255,  1067913)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2108", 56000, 2000000)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2109", 80000, 2000000)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2110", 80000, 2000000)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2111", 80000, 2000000)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2112", 80000, 2000000)
    ,@("CC", "1143405613", "AMAURY GUZMAN ACEVEDO",            "2112", 13867, 1300000)
    ,@("CC", "9148291",    "CARLOS ANDRES MARQUEZ MELENDEZ",   "2201", 80000, 2000000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row++
}

# Summary block above the table
$ws.Range("E11").Value = 549377   # VALOR MORA total
$ws.Range("C13").Value = 9        # Cant. Trabajadores
$ws.Range("F13").Value = 11       # Cant. Periodos

Write-Host ("Done. UsedRange: {0}" -f $ws.UsedRange.Address())
